# Conserto do erro com o rotulo da coluna 2050 nas tabelas e
# retirada das linhas com total das tabelas.
#
# Sheets 1-5 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)",
# "Emissoes Totais (MtCO2eq)") all had a stray numeric value
# (748.2245541884627) in E1 instead of the intended text label for the
# last column (2050, or 2041-2050 on the "Potencia Incremental" sheet).
# Sheets 1-4 additionally had a trailing "Total" row that should be removed.
# Sheet 6 ("Custo Total") also had a trailing "Total" row to remove.

$wb = $excel.ActiveWorkbook

# Map of sheet index -> corrected E1 label.
$labels = @{
    1 = "2050"
    2 = "2050"
    3 = "2050"
    4 = "2041-2050"
    5 = "2050"
}

foreach ($idx in $labels.Keys) {
    $ws = $wb.Worksheets.Item($idx)
    # A leading apostrophe forces Excel to store the value as text even
    # though it looks like a number, matching the intended text label
    # instead of the stray numeric value that was there before.
    $ws.Range("E1").Value = "'" + $labels[$idx]
}

# Remove the trailing "Total" row (row 13) from the first four sheets.
foreach ($idx in 1..4) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
}

# Remove the trailing "Total" row (row 4) from the "Custo Total" sheet.
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
